$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing row 51: E51 value correction (1.1787 -> 4.1787) ---
$ws.Range("E51").Value = 4.1787000000000001

# --- Row 53: fill in the remaining measurement columns + comment ---
$ws.Range("D53").Value = 534.47
$ws.Range("E53").Value = 4.1948
$ws.Range("F53").Value = 1198
$ws.Range("G53").Value = 175.7
$ws.Range("H53").Value = 2242
$ws.Range("I53").Value = "pos consistent as vert, focus tiny off, pmtB signals worse for more ypos when horz?, still not much junk or clustering"

# --- Row 54: new horz run ---
$ws.Range("A54").Value = "D20151105T212301"
$ws.Range("B54").Value = 5
$ws.Range("C54").Value = "H"
$ws.Range("D54").Value = 560.29999999999995
$ws.Range("E54").Value = 4.2214
$ws.Range("F54").Value = 1198
$ws.Range("G54").Value = 184.58
$ws.Range("H54").Value = 2366
$ws.Range("I54").Value = "same as last"
$ws.Range("J54").Value = "9um beads, use all signals"

# --- Row 55: new horz run ---
$ws.Range("A55").Value = "D20151105T214520"
$ws.Range("B55").Value = 5
$ws.Range("C55").Value = "H"
$ws.Range("D55").Value = 524.5
$ws.Range("E55").Value = 4.2076000000000002
$ws.Range("F55").Value = 1198
$ws.Range("G55").Value = 171.6
$ws.Range("H55").Value = 2207
$ws.Range("I55").Value = "all files look good, slightly lower conc, some a tiny bit out of focus"

# --- Move the active selection to the new bottom-most data entry row ---
$ws.Range("D56").Select()
